{"js": "// Update the \"Unidad 1: Algoritmos\" heading to \"Unidad 3: Diagramaci\u00f3n de algoritmos\"\n// The line is made of two runs: a bold \"Unidad 1:\" run and a plain \" Algoritmos\" run.\n// We replace each run's text in place (via search) so the existing run formatting\n// (bold vs. non-bold) is preserved exactly as in the original document.\n\nconst boldResults = context.document.body.search(\"Unidad 1:\", { matchCase: true });\nboldResults.load(\"items\");\nawait context.sync();\n\nif (boldResults.items.length > 0) {\n  boldResults.items[0].insertText(\"Unidad 3:\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst plainResults = context.document.body.search(\" Algoritmos\", { matchCase: true });\nplainResults.load(\"items\");\nawait context.sync();\n\nif (plainResults.items.length > 0) {\n  plainResults.items[0].insertText(\" Diagramaci\u00f3n de algoritmos\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the \"Unidad 1: Algoritmos\" heading to \"Unidad 3: Diagramaci\u00f3n de algoritmos\"\n# The line is made of two runs: a bold \"Unidad 1:\" run and a plain \" Algoritmos\" run.\n# Locate each piece of text with Find and overwrite just that range's text so the\n# surrounding run formatting (bold vs. non-bold) stays exactly as it was.\n\n$d = $word.ActiveDocument\n\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Text = \"Unidad 1:\"\n$find1.Forward = $true\n$find1.Wrap = 0\nif ($find1.Execute()) {\n    $rng1.Text = \"Unidad 3:\"\n}\n\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \" Algoritmos\"\n$find2.Forward = $true\n$find2.Wrap = 0\nif ($find2.Execute()) {\n    $rng2.Text = \" Diagramaci\u00f3n de algoritmos\"\n}\n"}
